$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update password column for all data rows from "Linda12345." to "Choucair12345."
$ws.Range("Q2").Value = "Choucair12345."
$ws.Range("Q3").Value = "Choucair12345."
$ws.Range("Q4").Value = "Choucair12345."

# Update apellido (surname) for Linda from "Perilla" to "Fernandez"
$ws.Range("C2").Value = "Fernandez"

# Move the active selection on the "data" sheet from A5 to D2
$ws.Activate()
$ws.Range("D2").Select()
